# Aggiorna Excel e HTML
# Adds a new "Accessory spleen" entry to the Spleen section of the glossary.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the new data row (row 19) ------------------------------------
$ws.Range("A19").Value = "Spleen"
$ws.Range("B19").Value = "Accessory spleen"
$ws.Range("C19").Value = "Clip 1 B-mode"
$ws.Range("D19").Value = "https://youtu.be/_FckFwJwynI "

# --- Hyperlink the new YouTube link cell -----------------------------------
$ws.Hyperlinks.Add($ws.Range("D19"), "https://youtu.be/_FckFwJwynI")
# Hyperlinks.Add applies its own direct formatting; reapply the workbook's
# named hyperlink style so the cell matches the rest of column D.
$ws.Range("D19").Style = "Collegamento ipertestuale"

# --- Update the active selection/scroll position left by the author -------
$ws.Activate()
$ws.Range("D23").Select()
